# edit.ps1 - Word COM-interop script implementing the SRS.docx revision:
#   * Bump the document date from 27-02-2024 to 1-04-2024
#   * Swap the Flask tech-stack bullet for NodeJS (bold label + new description)
#   * Rename the "Backend (Flask)" architecture bullet to "Backend (Node JS)"
#   * Replace the generic cloud-hosting sentence with a concrete Netlify.com mention
#
# Commit message: "Added and test Backend APIs"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Version/date line: "Version 1.0 | Date: 27-02-2024" -> "...: 1-04-2024"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute("27-02-2024", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $r.Text = "1-04-2024"
} else {
    Write-Output "WARNING: date string not found"
}

# ---------------------------------------------------------------------------
# 2) Tech-stack bullet: "Flask: A micro web framework for Python." ->
#    "NodeJS: A JavaScript library for building backend servers."
#    (keep the bold label / regular description run split)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute("Flask:", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $r.Text = "NodeJS:"
} else {
    Write-Output "WARNING: 'Flask:' not found"
}

$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute("A micro web framework for Python.", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $r.Text = "A JavaScript library for building backend servers."
} else {
    Write-Output "WARNING: Flask description not found"
}

# ---------------------------------------------------------------------------
# 3) Architecture bullet: "Backend (Flask)" -> "Backend (Node JS)"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute("Backend (Flask)", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $r.Text = "Backend (Node JS)"
} else {
    Write-Output "WARNING: 'Backend (Flask)' not found"
}

# ---------------------------------------------------------------------------
# 4) Hosting sentence: "The system will be hosted on a cloud platform
#    (e.g., AWS, Azure)." -> "The system will be hosted on a Netlify.com"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute("cloud platform (e.g., AWS, Azure).", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $r.Text = "Netlify.com"
} else {
    Write-Output "WARNING: hosting sentence not found"
}
